$d = $word.ActiveDocument

# The last paragraph of the document reads "Order ID: {d.paymentInfo.orderId"
# but is missing its closing "}". The sibling placeholder just above it
# ("Payment Method: {d.paymentInfo.cardType}") already closes its brace with
# a dedicated run, so we mirror that: add a new run containing just "}"
# immediately after the "orderId" text, using the exact same character
# formatting as the run that holds "orderId" (BC Sans / A6A6A6 / 8pt).

$rng = $d.Content
$found = $rng.Find.Execute("orderId", $false, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find 'orderId' text to anchor the new run on."
}

# Keep a handle on the formatted "orderId" run so we can clone its exact
# run properties (rFonts/color/sz/szCs) onto the new run.
$src = $rng.Duplicate
$insertStart = $rng.End
$srcLen = $src.End - $src.Start

# Collapse to the end of "orderId" and clone the formatted text there -
# this carries over the full run formatting (including complex-script
# sz/szCs) that setting individual Font properties does not reliably emit.
$rng.Collapse(0)
$rng.FormattedText = $src.FormattedText

# The clone above duplicated the source text ("orderId") verbatim; replace
# just the newly inserted copy with the closing brace, keeping formatting.
$newRun = $d.Range($insertStart, $insertStart + $srcLen)
$newRun.Text = "}"
